$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column style (column A) from the last existing row (366)
# down to the new rows (367:374) so the new cells match the existing
# formatting (centered, bordered, bold, date-time number format) used
# throughout column A.
$ws.Range("A366").Copy()
$ws.Range("A367:A374").PasteSpecial(-4122)

# New daily data rows covering through 2021-09-09 (serial date 44448),
# as per "aggiornamento a 9/09 compreso".
$data = @(
    @(44441, 6, 14, 124.0145274160687),
    @(44442, 0, 11, 97.43998582691115),
    @(44443, 0, 10, 88.58180529719195),
    @(44444, 0, 9, 79.72362476747276),
    @(44445, 0, 7, 62.00726370803437),
    @(44446, 0, 6, 53.14908317831517),
    @(44447, 3, 9, 79.72362476747276),
    @(44448, 5, 8, 70.86544423775356)
)

$r = 367
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
